$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8917768001556396
$ws.Range("B1").Value = 1.570632338523865
$ws.Range("C1").Value = 5.118063449859619
$ws.Range("D1").Value = 4.781161308288574
$ws.Range("E1").Value = 1.481326103210449
